$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.078.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.74%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.605.15'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.39%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.13'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.52%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -0.61%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.602.95'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.45%  '
$ws.Range("E10").Value = '  +3.15%  '
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("E12").Value = '  -1.72%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.343'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.19'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.076.83'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("E16").Value = '  -2.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.060.23'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.604.00'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.66%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '367.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.17%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.04'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.39'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.28'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.83'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.82%  '
$ws.Range("E24").Value = '  -3.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.27'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.81%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.735.66'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '583.83'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.32%  '
$ws.Range("E30").Value = '  -0.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0986'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.57%  '
$ws.Range("E32").Value = '  -5.38%  '
$ws.Range("E33").Value = '  -3.55%  '
$ws.Range("E34").Value = '  -3.21%  '
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("E36").Value = '  -4.54%  '
$ws.Range("E37").Value = '  -2.96%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '156.16'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.01'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.16%  '
$ws.Range("E40").Value = '  -1.63%  '
$ws.Range("E41").Value = '  +0.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.22'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.43%  '
$ws.Range("E43").Value = '  -4.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.09'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.13%  '
$ws.Range("E45").Value = '  -0.04%  '
$ws.Range("E46").Value = '  -2.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₆0285'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.71'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.55%  '
$ws.Range("E49").Value = '  -3.15%  '
$ws.Range("E50").Value = '  -1.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.41'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.75%  '
